$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '66.856.50'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  +3.31%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'" + '3.840.81'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  +5.04%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  +0.28%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'" + '422.23'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  +3.71%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'" + '128.65'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  -3.62%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'" + '3.836.66'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  +5.10%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'" + '0.608'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  -2.22%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D10').Value = "'" + '0.719'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  -1.44%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'" + '0.156'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  -4.62%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'" + '0.0000334'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  +0.45%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'" + '40.77'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  -3.10%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'" + '10.37'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  +3.89%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'" + '4.459.43'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  +5.00%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'" + '15.64'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  +15.53%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'" + '3.841.26'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  +4.59%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'" + '0.137'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  -0.63%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'" + '19.81'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  -1.26%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'" + '67.275.96'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  +3.80%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  -0.79%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'" + '408.19'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  -3.45%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'" + '14.88'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  -2.83%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'" + '84.05'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  -2.37%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'" + '3.04'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'" + '  +0.96%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'" + '37.32'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  +4.01%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'" + '10.03'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  +5.85%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'" + '3.23'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  +0.84%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'" + '9.55'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  +36.77%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'" + '5.44'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  +5.91%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'" + '739.73'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  +8.66%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'" + '13.19'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  +3.17%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('B33').Value = "'" + 'Hedera'
$ws.Range('B33').Style = 'Normal'
$ws.Range('C33').Value = "'" + 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('C33').Style = 'Normal'
$ws.Range('D33').Value = "'" + '0.121'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  +3.22%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('B34').Value = "'" + 'Toncoin'
$ws.Range('B34').Style = 'Normal'
$ws.Range('C34').Value = "'" + 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('C34').Style = 'Normal'
$ws.Range('D34').Value = "'" + '2.69'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  -1.13%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  -0.11%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  -5.94%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'" + '38.46'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  -7.72%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('B38').Value = "'" + 'NEARProtocol'
$ws.Range('B38').Style = 'Normal'
$ws.Range('C38').Value = "'" + 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('C38').Style = 'Normal'
$ws.Range('D38').Value = "'" + '5.49'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  +23.28%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('B39').Value = "'" + 'OKB'
$ws.Range('B39').Style = 'Normal'
$ws.Range('C39').Value = "'" + 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('C39').Style = 'Normal'
$ws.Range('D39').Value = "'" + '55.55'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  -0.69%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'" + '0.0₃0729'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  +9.51%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  -2.51%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'" + '2.89'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  -2.19%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  +0.79%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'" + '3.36'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  +0.63%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  -4.95%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'" + '0.317'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  +8.44%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'" + '3.10'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  -0.17%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'" + '2.04'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  -2.53%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'" + '140.96'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'" + '  -2.16%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'" + '2.81'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  -0.54%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').Value = "'" + 'WEMIXToken'
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = "'" + 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = "'" + '2.54'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  +0.70%  '
$ws.Range('E51').Style = 'Normal'
